$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update amount values
$ws.Range("D2").Value = 11
$ws.Range("D4").Value = 8
$ws.Range("D5").Value = 9

# Update size and color for the last row (shirt)
$ws.Range("C7").Value = "s/m"
$ws.Range("E7").Value = "red"
